# Updates cryptos list values (Price and Volume(1h) columns) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '71.208.39'
$ws.Range("E2").Value = '  +0.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.571.42'
$ws.Range("E3").Value = '  -0.50%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.57'
$ws.Range("E5").Value = '  +0.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.69'
$ws.Range("E6").Value = '  +0.83%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.519'
$ws.Range("E8").Value = '  +1.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.571.22'
$ws.Range("E9").Value = '  -0.54%  '
$ws.Range("E10").Value = '  +0.28%  '
$ws.Range("E11").Value = '  -0.70%  '
$ws.Range("E12").Value = '  +1.30%  '
$ws.Range("E13").Value = '  +1.62%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.035.89'
$ws.Range("E14").Value = '  -0.89%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '70.868.65'
$ws.Range("E15").Value = '  -0.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000180'
$ws.Range("E16").Value = '  -2.25%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.54'
$ws.Range("E17").Value = '  +1.06%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.578.32'
$ws.Range("E18").Value = '  -0.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.95'
$ws.Range("E19").Value = '  +3.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.51'
$ws.Range("E20").Value = '  -2.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '357.89'
$ws.Range("E21").Value = '  -2.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.96'
$ws.Range("E22").Value = '  -0.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.09'
$ws.Range("E23").Value = '  +4.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.01'
$ws.Range("E24").Value = '  +0.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.50'
$ws.Range("E25").Value = '  -0.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.10'
$ws.Range("E26").Value = '  -1.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.25'
$ws.Range("E27").Value = '  -0.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.667.19'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.997'
$ws.Range("E29").Value = '  -0.67%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0929'
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.00'
$ws.Range("E31").Value = '  +1.78%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '478.00'
$ws.Range("E32").Value = '  -0.81%  '
$ws.Range("E33").Value = '  -1.46%  '
$ws.Range("E34").Value = '  +0.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.19%  '
$ws.Range("E36").Value = '  +4.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '160.02'
$ws.Range("E37").Value = '  +0.98%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.11'
$ws.Range("E38").Value = '  +1.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.77'
$ws.Range("E39").Value = '  -0.53%  '
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.93'
$ws.Range("E41").Value = '  +3.32%  '
$ws.Range("E42").Value = '  +0.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.64'
$ws.Range("E43").Value = '  -3.96%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.39'
$ws.Range("E44").Value = '  -4.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.18'
$ws.Range("E45").Value = '  -11.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.59'
$ws.Range("E46").Value = '  -0.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '145.87'
$ws.Range("E47").Value = '  -1.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.542'
$ws.Range("E48").Value = '  +1.79%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.57'
$ws.Range("E49").Value = '  -0.61%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.64'
$ws.Range("E50").Value = '  -0.54%  '
$ws.Range("E51").Value = '  +0.71%  '
